# Update "Case and Fatality Demographics Data" workbook (2021-12-03 refresh).
# Only the Fatalities sheets carry new underlying data; the Cases sheets are
# untouched numerically (only cosmetic/view metadata differs upstream, which
# is not meaningful to reproduce via the object model).

$wb = $excel.ActiveWorkbook

# --- Fatalities by Age Group ---
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B3").Value  = 22
$ws.Range("B4").Value  = 79
$ws.Range("B5").Value  = 662
$ws.Range("B6").Value  = 2130
$ws.Range("B7").Value  = 5164
$ws.Range("B8").Value  = 9824
$ws.Range("B9").Value  = 7481
$ws.Range("B10").Value = 8775
$ws.Range("B11").Value = 9290
$ws.Range("B12").Value = 8784
$ws.Range("B13").Value = 20657
$ws.Range("B15").Value = 72882

# --- Fatalities by Gender ---
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 30520
$ws.Range("B3").Value = 42361
$ws.Range("B5").Value = 72882

# --- Fatalities by Race-Ethnicity ---
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 1345
$ws.Range("B3").Value = 7725
$ws.Range("B4").Value = 31810
$ws.Range("B5").Value = 433
$ws.Range("B6").Value = 31524
